$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at C ("FLC"), shifting NP,NC,NR,NFRec,NCRec,NIFPar,NRRec,NRNRec,Final_Score
# one column to the right (D..L).
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "FLC"

# New FLC feature values for the data rows.
$ws.Cells.Item(2, 3).Value = 15.75
$ws.Cells.Item(3, 3).Value = 18.25
$ws.Cells.Item(4, 3).Value = 18.25

# The Final_Score column (now column L) text needs updating for rows 3 and 4
# because the score now also accounts for the new feature. Force these as
# plain text (not a parsed percentage number) and keep the default style,
# matching the other Final_Score cells.
$c = $ws.Cells.Item(3, 12)
$c.NumberFormat = "@"
$c.Value = "45.49%"
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 12)
$c.NumberFormat = "@"
$c.Value = "75.91%"
$c.Style = "Normal"
